$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for the columns that change (D, J, K, L, M, P)
# before overwriting, since the edit permutes data across rows 2-46.
$cols = @("D","J","K","L","M","P")
$orig = @{}
for ($r = 2; $r -le 46; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

# Map: destination row -> source row (values to copy from)
$rowMap = @{
    2 = 4
    3 = 25
    4 = 19
    5 = 27
    6 = 46
    7 = 41
    8 = 16
    9 = 23
    10 = 32
    11 = 6
    12 = 38
    13 = 39
    14 = 33
    15 = 8
    16 = 45
    17 = 13
    18 = 22
    19 = 26
    20 = 10
    21 = 31
    22 = 15
    23 = 30
    24 = 29
    25 = 37
    26 = 7
    27 = 3
    28 = 40
    29 = 2
    30 = 9
    31 = 36
    32 = 43
    33 = 11
    34 = 35
    35 = 24
    36 = 28
    37 = 18
    38 = 14
    39 = 12
    40 = 20
    41 = 21
    42 = 34
    43 = 17
    44 = 42
    45 = 44
    46 = 5
}

foreach ($dest in $rowMap.Keys) {
    $src = $rowMap[$dest]
    $srcVals = $orig[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value2 = $srcVals[$c]
    }
}
